# Delete the "Crime against women" column (column I) from Sheet1.
# This shifts the trailing "Safety score" column (J) left into column I,
# matching the target workbook layout (A:I instead of A:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(9).Delete()

# Select the freshly-shifted column, mirroring how Excel leaves the
# selection after deleting an entire column (whole-column sqref).
$ws.Columns.Item(9).Select()
